$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of the price table: for each data row (2-21), update
# Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) with the new values.

$ws.Cells.Item(2, 4).Value = 44875
$ws.Cells.Item(2, 13).Value = 400
$ws.Cells.Item(2, 14).Value = 7000
$ws.Cells.Item(2, 15).Value = 7500
$ws.Cells.Item(2, 16).Value = 7250
$ws.Cells.Item(2, 19).Value = 3625

$ws.Cells.Item(3, 4).Value = 44482
$ws.Cells.Item(3, 13).Value = 240
$ws.Cells.Item(3, 14).Value = 10000
$ws.Cells.Item(3, 15).Value = 11000
$ws.Cells.Item(3, 16).Value = 10500
$ws.Cells.Item(3, 19).Value = 5250

$ws.Cells.Item(4, 4).Value = 44489
$ws.Cells.Item(4, 13).Value = 160
$ws.Cells.Item(4, 14).Value = 9500
$ws.Cells.Item(4, 15).Value = 10000
$ws.Cells.Item(4, 16).Value = 9750
$ws.Cells.Item(4, 19).Value = 4875

$ws.Cells.Item(5, 4).Value = 45211
$ws.Cells.Item(5, 13).Value = 200
$ws.Cells.Item(5, 14).Value = 10000
$ws.Cells.Item(5, 15).Value = 11000
$ws.Cells.Item(5, 16).Value = 10500
$ws.Cells.Item(5, 19).Value = 5250

$ws.Cells.Item(6, 4).Value = 44881
$ws.Cells.Item(6, 13).Value = 440
$ws.Cells.Item(6, 14).Value = 6000
$ws.Cells.Item(6, 15).Value = 7000
$ws.Cells.Item(6, 16).Value = 6500
$ws.Cells.Item(6, 19).Value = 3250

$ws.Cells.Item(7, 4).Value = 44461
$ws.Cells.Item(7, 13).Value = 200
$ws.Cells.Item(7, 14).Value = 11000
$ws.Cells.Item(7, 15).Value = 12000
$ws.Cells.Item(7, 16).Value = 11500
$ws.Cells.Item(7, 19).Value = 5750

$ws.Cells.Item(8, 4).Value = 44497
$ws.Cells.Item(8, 13).Value = 500
$ws.Cells.Item(8, 14).Value = 9000
$ws.Cells.Item(8, 15).Value = 10000
$ws.Cells.Item(8, 16).Value = 9500
$ws.Cells.Item(8, 19).Value = 4750

$ws.Cells.Item(9, 4).Value = 44895
$ws.Cells.Item(9, 13).Value = 240
$ws.Cells.Item(9, 14).Value = 3000
$ws.Cells.Item(9, 15).Value = 3500
$ws.Cells.Item(9, 16).Value = 3250
$ws.Cells.Item(9, 19).Value = 1625

$ws.Cells.Item(10, 4).Value = 45204
$ws.Cells.Item(10, 13).Value = 400
$ws.Cells.Item(10, 14).Value = 9000
$ws.Cells.Item(10, 15).Value = 10000
$ws.Cells.Item(10, 16).Value = 9500
$ws.Cells.Item(10, 19).Value = 4750

$ws.Cells.Item(11, 4).Value = 44889
$ws.Cells.Item(11, 13).Value = 460
$ws.Cells.Item(11, 14).Value = 3500
$ws.Cells.Item(11, 15).Value = 4000
$ws.Cells.Item(11, 16).Value = 3750
$ws.Cells.Item(11, 19).Value = 1875

$ws.Cells.Item(12, 4).Value = 44454
$ws.Cells.Item(12, 13).Value = 160
$ws.Cells.Item(12, 14).Value = 12000
$ws.Cells.Item(12, 15).Value = 13000
$ws.Cells.Item(12, 16).Value = 12500
$ws.Cells.Item(12, 19).Value = 6250

$ws.Cells.Item(13, 4).Value = 44882
$ws.Cells.Item(13, 13).Value = 440
$ws.Cells.Item(13, 14).Value = 6000
$ws.Cells.Item(13, 15).Value = 7000
$ws.Cells.Item(13, 16).Value = 6500
$ws.Cells.Item(13, 19).Value = 3250

$ws.Cells.Item(14, 4).Value = 44517
$ws.Cells.Item(14, 13).Value = 400
$ws.Cells.Item(14, 14).Value = 5500
$ws.Cells.Item(14, 15).Value = 6000
$ws.Cells.Item(14, 16).Value = 5750
$ws.Cells.Item(14, 19).Value = 2875

$ws.Cells.Item(15, 4).Value = 44455
$ws.Cells.Item(15, 13).Value = 200
$ws.Cells.Item(15, 14).Value = 12000
$ws.Cells.Item(15, 15).Value = 13000
$ws.Cells.Item(15, 16).Value = 12500
$ws.Cells.Item(15, 19).Value = 6250

$ws.Cells.Item(16, 4).Value = 44475
$ws.Cells.Item(16, 13).Value = 240
$ws.Cells.Item(16, 14).Value = 11000
$ws.Cells.Item(16, 15).Value = 12000
$ws.Cells.Item(16, 16).Value = 11500
$ws.Cells.Item(16, 19).Value = 5750

$ws.Cells.Item(17, 4).Value = 44819
$ws.Cells.Item(17, 13).Value = 240
$ws.Cells.Item(17, 14).Value = 11000
$ws.Cells.Item(17, 15).Value = 12000
$ws.Cells.Item(17, 16).Value = 11500
$ws.Cells.Item(17, 19).Value = 5750

$ws.Cells.Item(18, 4).Value = 45203
$ws.Cells.Item(18, 13).Value = 1000
$ws.Cells.Item(18, 14).Value = 9000
$ws.Cells.Item(18, 15).Value = 10000
$ws.Cells.Item(18, 16).Value = 9500
$ws.Cells.Item(18, 19).Value = 4750

$ws.Cells.Item(19, 4).Value = 44490
$ws.Cells.Item(19, 13).Value = 400
$ws.Cells.Item(19, 14).Value = 9500
$ws.Cells.Item(19, 15).Value = 10000
$ws.Cells.Item(19, 16).Value = 9750
$ws.Cells.Item(19, 19).Value = 4875

$ws.Cells.Item(20, 4).Value = 44874
$ws.Cells.Item(20, 13).Value = 300
$ws.Cells.Item(20, 14).Value = 7500
$ws.Cells.Item(20, 15).Value = 8000
$ws.Cells.Item(20, 16).Value = 7750
$ws.Cells.Item(20, 19).Value = 3875

$ws.Cells.Item(21, 4).Value = 44818
$ws.Cells.Item(21, 13).Value = 200
$ws.Cells.Item(21, 14).Value = 11000
$ws.Cells.Item(21, 15).Value = 12000
$ws.Cells.Item(21, 16).Value = 11500
$ws.Cells.Item(21, 19).Value = 5750
